$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and G contain numeric-looking text values in the source data.
# Force them to Text format first so Excel keeps them as exact strings
# instead of converting to floating-point numbers.
$ws.Range("D2,G2,G3,D4,G4,D5,G5,D6,G6,D7,G7,D8,G8,D9,G9,D10,G10,D11,G11,D12,G12,D13,G13,D14,G14,D15,G15,D16,G16,D17,G17,D18,G18,D19,G19,D20,G20,D21,G21,D22,G22,D23,G23,G24,D25,G25,D26,G26,D27,G27,G28,G29,G30,G31,G32,G33,G34,G35,G36,G37,G38,G39,D40,G40,D41,G41,D42,G42,D43,G43,D44,G44,D45,G45,D46,G46,G47,D48,G48,D49,G49,G50,G51").NumberFormat = "@"

$ws.Range("D2").Value = "260.95"
$ws.Range("F2").Value = "16-12-2022"
$ws.Range("G2").Value = "0"
$ws.Range("F3").Value = "16-12-2022"
$ws.Range("G3").Value = "0"
$ws.Range("D4").Value = "6.182"
$ws.Range("F4").Value = "16-12-2022"
$ws.Range("G4").Value = "0"
$ws.Range("D5").Value = "0.06093"
$ws.Range("F5").Value = "16-12-2022"
$ws.Range("G5").Value = "0"
$ws.Range("D6").Value = "6.749"
$ws.Range("F6").Value = "16-12-2022"
$ws.Range("G6").Value = "0"
$ws.Range("D7").Value = "3.442"
$ws.Range("F7").Value = "16-12-2022"
$ws.Range("G7").Value = "0"
$ws.Range("D8").Value = "1.353"
$ws.Range("F8").Value = "16-12-2022"
$ws.Range("G8").Value = "0"
$ws.Range("D9").Value = "0.7987"
$ws.Range("F9").Value = "16-12-2022"
$ws.Range("G9").Value = "0"
$ws.Range("D10").Value = "0.1579"
$ws.Range("F10").Value = "16-12-2022"
$ws.Range("G10").Value = "0"
$ws.Range("D11").Value = "0.08102"
$ws.Range("F11").Value = "16-12-2022"
$ws.Range("G11").Value = "0"
$ws.Range("D12").Value = "0.03365"
$ws.Range("F12").Value = "16-12-2022"
$ws.Range("G12").Value = "0"
$ws.Range("D13").Value = "0.03088"
$ws.Range("F13").Value = "16-12-2022"
$ws.Range("G13").Value = "0"
$ws.Range("D14").Value = "0.09319"
$ws.Range("F14").Value = "16-12-2022"
$ws.Range("G14").Value = "0"
$ws.Range("D15").Value = "3.855"
$ws.Range("F15").Value = "16-12-2022"
$ws.Range("G15").Value = "0"
$ws.Range("D16").Value = "0.001700"
$ws.Range("F16").Value = "16-12-2022"
$ws.Range("G16").Value = "0"
$ws.Range("D17").Value = "0.04834"
$ws.Range("F17").Value = "16-12-2022"
$ws.Range("G17").Value = "0"
$ws.Range("D18").Value = "0.0006149"
$ws.Range("F18").Value = "16-12-2022"
$ws.Range("G18").Value = "0"
$ws.Range("D19").Value = "0.006215"
$ws.Range("F19").Value = "16-12-2022"
$ws.Range("G19").Value = "0"
$ws.Range("D20").Value = "0.001097"
$ws.Range("F20").Value = "16-12-2022"
$ws.Range("G20").Value = "0"
$ws.Range("D21").Value = "0.003399"
$ws.Range("F21").Value = "16-12-2022"
$ws.Range("G21").Value = "0"
$ws.Range("D22").Value = "0.0001500"
$ws.Range("F22").Value = "16-12-2022"
$ws.Range("G22").Value = "0"
$ws.Range("D23").Value = "3.682"
$ws.Range("F23").Value = "16-12-2022"
$ws.Range("G23").Value = "0"
$ws.Range("F24").Value = "16-12-2022"
$ws.Range("G24").Value = "0"
$ws.Range("D25").Value = "0.3361"
$ws.Range("F25").Value = "16-12-2022"
$ws.Range("G25").Value = "0"
$ws.Range("D26").Value = "0.1260"
$ws.Range("F26").Value = "16-12-2022"
$ws.Range("G26").Value = "0"
$ws.Range("D27").Value = "0.0006110"
$ws.Range("F27").Value = "16-12-2022"
$ws.Range("G27").Value = "0"
$ws.Range("F28").Value = "16-12-2022"
$ws.Range("G28").Value = "0"
$ws.Range("F29").Value = "16-12-2022"
$ws.Range("G29").Value = "0"
$ws.Range("F30").Value = "16-12-2022"
$ws.Range("G30").Value = "0"
$ws.Range("F31").Value = "16-12-2022"
$ws.Range("G31").Value = "0"
$ws.Range("F32").Value = "16-12-2022"
$ws.Range("G32").Value = "0"
$ws.Range("F33").Value = "16-12-2022"
$ws.Range("G33").Value = "0"
$ws.Range("F34").Value = "16-12-2022"
$ws.Range("G34").Value = "0"
$ws.Range("F35").Value = "16-12-2022"
$ws.Range("G35").Value = "0"
$ws.Range("F36").Value = "16-12-2022"
$ws.Range("G36").Value = "0"
$ws.Range("F37").Value = "16-12-2022"
$ws.Range("G37").Value = "0"
$ws.Range("F38").Value = "16-12-2022"
$ws.Range("G38").Value = "0"
$ws.Range("F39").Value = "16-12-2022"
$ws.Range("G39").Value = "0"
$ws.Range("D40").Value = "0.04577"
$ws.Range("F40").Value = "16-12-2022"
$ws.Range("G40").Value = "0"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.007159"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("F41").Value = "16-12-2022"
$ws.Range("G41").Value = "0"
$ws.Range("D42").Value = "0.003899"
$ws.Range("F42").Value = "16-12-2022"
$ws.Range("G42").Value = "0"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "0.1119"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("F43").Value = "16-12-2022"
$ws.Range("G43").Value = "0"
$ws.Range("D44").Value = "0.01008"
$ws.Range("F44").Value = "16-12-2022"
$ws.Range("G44").Value = "0"
$ws.Range("D45").Value = "0.002969"
$ws.Range("F45").Value = "16-12-2022"
$ws.Range("G45").Value = "0"
$ws.Range("D46").Value = "0.00005923"
$ws.Range("F46").Value = "16-12-2022"
$ws.Range("G46").Value = "0"
$ws.Range("F47").Value = "16-12-2022"
$ws.Range("G47").Value = "0"
$ws.Range("D48").Value = "0.6999"
$ws.Range("F48").Value = "16-12-2022"
$ws.Range("G48").Value = "0"
$ws.Range("D49").Value = "0.07337"
$ws.Range("F49").Value = "16-12-2022"
$ws.Range("G49").Value = "0"
$ws.Range("F50").Value = "16-12-2022"
$ws.Range("G50").Value = "0"
$ws.Range("F51").Value = "16-12-2022"
$ws.Range("G51").Value = "0"
